# MustafaSheikh_Resume.docx edit:
#  1. Trim "simultaneous" from the HIL bullet's 3-platform sentence.
#  2. Turn off "overflow punctuation" (w:overflowPunct) on the three
#     paragraph styles that still had it on (Normal, TOC Heading, No
#     Spacing) — this is what produced the left-justified Skills layout.

$d = $word.ActiveDocument

# 1) Shrink the bullet text.
$d.Content.Find.Execute(
    "Slashed HIL part costs by 75% through BOM management of 3 simultaneous vehicle platforms, involving 100s of components.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Slashed HIL part costs by 75% through BOM management of 3 vehicle platforms, involving 100s of components.",
    2
) | Out-Null

# 2) Flip w:overflowPunct true -> false on the affected styles.
$styleNames = @("Normal", "TOC Heading", "No Spacing")
foreach ($styleName in $styleNames) {
    $style = $d.Styles($styleName)
    $style.ParagraphFormat.HangingPunctuation = $false
}
